# Fix Results.xlsx to match modified H-Score thresholds.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("H-Score")

# Row 4
$ws.Range("C4").Value = 814
$ws.Range("D4").Value = 56
$ws.Range("E4").Value = 18
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 0.907
$ws.Range("H4").Value = 0.062
$ws.Range("I4").Value = 0.02
$ws.Range("J4").Value = 0.01
$ws.Range("K4").Value = 13

# Row 5
$ws.Range("C5").Value = 566
$ws.Range("D5").Value = 169
$ws.Range("E5").Value = 51
$ws.Range("F5").Value = 49
$ws.Range("G5").Value = 0.678
$ws.Range("H5").Value = 0.202
$ws.Range("I5").Value = 0.061
$ws.Range("J5").Value = 0.059
$ws.Range("K5").Value = 50

# Row 6
$ws.Range("C6").Value = 1380
$ws.Range("D6").Value = 225
$ws.Range("E6").Value = 69
$ws.Range("F6").Value = 58
$ws.Range("G6").Value = 0.797
$ws.Range("H6").Value = 0.13
$ws.Range("I6").Value = 0.04
$ws.Range("J6").Value = 0.033
$ws.Range("K6").Value = 31

# Row 7
$ws.Range("C7").Value = 414
$ws.Range("D7").Value = 109
$ws.Range("E7").Value = 56
$ws.Range("F7").Value = 55
$ws.Range("G7").Value = 0.653
$ws.Range("H7").Value = 0.172
$ws.Range("I7").Value = 0.088
$ws.Range("J7").Value = 0.087
$ws.Range("K7").Value = 61

# Row 8
$ws.Range("C8").Value = 843
$ws.Range("D8").Value = 263
$ws.Range("E8").Value = 72
$ws.Range("F8").Value = 102
$ws.Range("G8").Value = 0.659
$ws.Range("H8").Value = 0.205
$ws.Range("I8").Value = 0.056
$ws.Range("J8").Value = 0.08
$ws.Range("K8").Value = 56

# Row 9
$ws.Range("C9").Value = 1257
$ws.Range("D9").Value = 372
$ws.Range("E9").Value = 128
$ws.Range("F9").Value = 157
$ws.Range("G9").Value = 0.657
$ws.Range("H9").Value = 0.194
$ws.Range("I9").Value = 0.067
$ws.Range("J9").Value = 0.082
$ws.Range("K9").Value = 57

# Row 10
$ws.Range("C10").Value = 527
$ws.Range("D10").Value = 57
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 13
$ws.Range("G10").Value = 0.868
$ws.Range("H10").Value = 0.094
$ws.Range("I10").Value = 0.016
$ws.Range("J10").Value = 0.021
$ws.Range("K10").Value = 19

# Row 11
$ws.Range("C11").Value = 927
$ws.Range("D11").Value = 164
$ws.Range("E11").Value = 66
$ws.Range("F11").Value = 45
$ws.Range("G11").Value = 0.771
$ws.Range("H11").Value = 0.136
$ws.Range("I11").Value = 0.055
$ws.Range("J11").Value = 0.037
$ws.Range("K11").Value = 36

# Row 12
$ws.Range("C12").Value = 1454
$ws.Range("D12").Value = 221
$ws.Range("E12").Value = 76
$ws.Range("F12").Value = 58
$ws.Range("G12").Value = 0.804
$ws.Range("H12").Value = 0.122
$ws.Range("I12").Value = 0.042
$ws.Range("J12").Value = 0.032
$ws.Range("K12").Value = 30
